$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-6 (columns A:S), replacing/extending the existing single
# data row with 5 rows of trade records.
$data = @(
    @("AUR_3.0_Call",  "SELL", 45875.53912122138, 12, 3.2,                3.094540909255746,  3,  6.07,  1, 25, 0.1054590907442545,  100, 45876.29025811746, 6,  3.2,                3.266588328671542,  0, 0, -6),
    @("AUR_14.0_Call", "SELL", 45875.53912122138, 12, 0.05,               0.02877876720203026, 14, 6.07,  1, 25, 0.02122123279796975, 100, 45876.4170691109,  10, 0.05,               0.1517281683239983, 0, 0, -2),
    @("AUR_3.0_Call",  "SELL", 45876.33251067522, 7,  3.2,                3.123114092017505,  3,  6.1,   1, 25, 0.07688590798249484, 100, 45876.4170691109,  10, 3.2,                3.20566315489788,   0, 0, 3),
    @("AUR_11.0_Call", "SELL", 45876.4593148956,  11, 0.07000000000000001, 0.03679925764654384, 11, 6.085, 1, 25, 0.03320074235345617, 100, 45876.50165914217, 12, 0.07000000000000001, 0.08659202483859844, 0, 0, 1),
    @("AUR_14.0_Call", "SELL", 45876.4593148956,  11, 0.05,               0.005796167528757829, 14, 6.085, 1, 25, 0.04420383247124218, 100, 45877.27325131699, 6,  0.05,               0.09228680544472716, 0, 0, -5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
    # Preserve the datetime-formatted style on entry_time (col C) and exit_time (col M)
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 13).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
